$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shape = $s.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text -eq "5 Button") {
                    $shape.TextFrame.TextRange.Text = "6 Button Toggle"
                }
            }
        }
    }
}
